# Add the new "2022-Q3" quarterly sheet right after "总计", shifting the
# existing quarter sheets down, and update the "总计" summary sheet with a
# new row for 2022-Q3.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)

# A cell that already carries the workbook's "header / row-number" look
# (bold font, thin border, centered) so we can clone its formatting onto
# the new sheet without inventing a brand-new style.
$styleDonor = $totalSheet.Range("B1")

# ----------------------------------------------------------------------
# 1. Insert a brand new worksheet right after "总计" (the first sheet) and
#    name it "2022-Q3". All the other quarter sheets shift right by one
#    position automatically.
# ----------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row for the new sheet (matches the layout used by the other
# quarterly fund-holding sheets).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$styleDonor.Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3Rows = @(
    @("900090", "中信卓越成长两年持有期混合B", "51.50", "93.14", "3.12", "1.6068", 9),
    @("166301", "华商新趋势优选灵活配置混合",   "57.69", "75.38", "1.98", "1.1423", 5),
    @("000390", "华商优势行业混合",             "26.43", "83.59", "2.06", "0.5445", 10),
    @("900010", "中信卓越成长两年持有期混合A", "14.24", "93.14", "3.12", "0.4443", 9),
    @("900100", "中信卓越成长两年持有期混合C", "4.61",  "93.14", "3.12", "0.1438", 9),
    @("008488", "华商恒益稳健混合",             "2.39",  "52.53", "1.61", "0.0385", 9),
    @("008629", "大成景瑞稳健配置混合A",       "0.89",  "21.44", "1.37", "0.0122", 7),
    @("008630", "大成景瑞稳健配置混合C",       "0.75",  "21.44", "1.37", "0.0103", 7),
    @("001231", "银华泰利灵活配置混合A",       "0.96",  "24.09", "0.77", "0.0074", 9),
    @("003063", "银华通利灵活配置混合C",       "0.30",  "26.56", "0.77", "0.0023", 10),
    @("003062", "银华通利灵活配置混合A",       "0.22",  "26.56", "0.77", "0.0017", 10),
    @("002328", "银华泰利灵活配置混合C",       "0.03",  "24.09", "0.77", "0.0002", 9)
)

# Columns B-G on the quarterly fund-holding sheets are stored as text
# (even the numeric-looking scale/position figures), so force a text
# number format before writing any values into them, then clear the
# number format back off afterwards so the cells end up with the
# workbook's default (no) style, exactly like the other sheets.
$q3.Range("B2:G13").NumberFormat = "@"

$r = 2
foreach ($row in $q3Rows) {
    $q3.Range("A$r").Value = ($r - 2)
    $q3.Range("B$r").Value = $row[0]
    $q3.Range("C$r").Value = $row[1]
    $q3.Range("D$r").Value = $row[2]
    $q3.Range("E$r").Value = $row[3]
    $q3.Range("F$r").Value = $row[4]
    $q3.Range("G$r").Value = $row[5]
    $q3.Range("H$r").Value = $row[6]
    $r = $r + 1
}

$q3.Range("B2:G13").ClearFormats()

$styleDonor.Copy()
$q3.Range("A2:A13").PasteSpecial(-4122)

# ----------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row right under the header for
#    2022-Q3 and push the rest of the rows down by one.
# ----------------------------------------------------------------------
$totalSheet.Range("A2:D2").Rows.Insert()

# Inserting a row clones formatting from the row above (the bold header
# row), which is not what any of the other data rows look like, so wipe
# it back to the sheet's default before writing the new values.
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 12
$totalSheet.Range("D2").Value = 3.95

# Column A always carries the bold/bordered "row index" style, just like
# every other row on this sheet.
$styleDonor.Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# The "#" column (A) holds literal sequential numbers rather than a
# formula, so renumber the rows that got pushed down by the insert.
for ($row = 3; $row -le 8; $row++) {
    $totalSheet.Range("A$row").Value = $row - 2
}
